# Apply the changes described by the commit:
# "import export tanah, gedung, kendaraan fix"
#
# Concretely, for this worksheet template:
#  - Column P (header cell P1, which only carried a fill style / no value)
#    is removed from the used range entirely.
#  - N2 changes from the text value "-" to the numeric value 0.
#  - The active selection / view moves to N3 (scrolled so column E is
#    the left-most visible column).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused P1 cell (style-only, no content) so the sheet's
# dimension/used range shrinks back from A1:P2 to A1:O2.
$ws.Range("P1").Clear()

# N2 used to hold the text "-" (shared string); it now holds a literal
# number 0.
$ws.Range("N2").Value = 0

# Update the view/selection to match the saved workbook state.
$win = $excel.ActiveWindow
$win.ScrollColumn = 5
$win.ScrollRow = 1
$ws.Range("N3").Select()
